$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 127.15385
$ws.Range("I38").Value = 83.583336
$ws.Range("J38").Value = 650
$ws.Range("K38").Value = 250.750008
$ws.Range("L38").Value = 1950
$ws.Range("M38").Value = 121.249992
$ws.Range("N38").Value = -2694
$ws.Range("H58").Value = 720.4545000000001
$ws.Range("I58").Value = 495.83334
$ws.Range("J58").Value = 990
$ws.Range("K58").Value = 1487.50002
$ws.Range("L58").Value = 2970
$ws.Range("M58").Value = -1337.50002
$ws.Range("N58").Value = -3270
$ws.Range("H87").Value = 29792
$ws.Range("J87").Value = 29792
$ws.Range("L87").Value = 29792
$ws.Range("N87").Value = -32288
$ws.Range("H90").Value = 29792
$ws.Range("J90").Value = 29792
$ws.Range("L90").Value = 89376
$ws.Range("N90").Value = -101856
$ws.Range("H98").Value = 3542.6843
$ws.Range("I98").Value = 3391.9062
$ws.Range("J98").Value = 4346.8335
$ws.Range("K98").Value = 3391.9062
$ws.Range("L98").Value = 4346.8335
$ws.Range("M98").Value = -1893.9062
$ws.Range("N98").Value = -7342.8335
$ws.Range("H122").Value = 3542.6843
$ws.Range("I122").Value = 3391.9062
$ws.Range("J122").Value = 4346.8335
$ws.Range("K122").Value = 10175.7186
$ws.Range("L122").Value = 13040.5005
$ws.Range("M122").Value = -7725.7186
$ws.Range("N122").Value = -17940.5005
$ws.Range("H125").Value = 1014.3929
$ws.Range("I125").Value = 689.1429000000001
$ws.Range("J125").Value = 1990.1428
$ws.Range("K125").Value = 6202.2861
$ws.Range("L125").Value = 17911.2852
$ws.Range("M125").Value = -3742.2861
$ws.Range("N125").Value = -22831.2852
$ws.Range("H137").Value = 18870236
$ws.Range("I137").Value = 1373.6666
$ws.Range("K137").Value = 4120.9998
$ws.Range("M137").Value = -1570.9998
$ws.Range("H138").Value = 3177.4658
$ws.Range("I138").Value = 2477.525
$ws.Range("J138").Value = 4025.879
$ws.Range("K138").Value = 7432.575000000001
$ws.Range("L138").Value = 12077.637
$ws.Range("M138").Value = -2292.575000000001
$ws.Range("N138").Value = -22357.637

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 15630991
$ws.Range("I74").Value = 25001250
$ws.Range("J74").Value = 13890.75
$ws.Range("K74").Value = 25001250
$ws.Range("L74").Value = 13890.75
$ws.Range("M74").Value = -25000376
$ws.Range("N74").Value = -15638.75
$ws.Range("H77").Value = 15630991
$ws.Range("I77").Value = 25001250
$ws.Range("J77").Value = 13890.75
$ws.Range("K77").Value = 125006250
$ws.Range("L77").Value = 69453.75
$ws.Range("M77").Value = -125001882
$ws.Range("N77").Value = -78189.75
$ws.Range("H132").Value = 1727161.8
$ws.Range("I132").Value = 2655156
$ws.Range("K132").Value = 7965468
$ws.Range("M132").Value = -7962938

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 343.33334
$ws.Range("I94").Value = 298.75
$ws.Range("J94").Value = 700
$ws.Range("K94").Value = 298.75
$ws.Range("L94").Value = 700
$ws.Range("M94").Value = 152.25
$ws.Range("N94").Value = -1602
$ws.Range("H107").Value = 347240.4
$ws.Range("I107").Value = 524724.75
$ws.Range("J107").Value = 4949.143
$ws.Range("K107").Value = 524724.75
$ws.Range("L107").Value = 4949.143
$ws.Range("M107").Value = -522804.75
$ws.Range("N107").Value = -8789.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4659.8184
$ws.Range("I31").Value = 3142.8
$ws.Range("J31").Value = 5924
$ws.Range("K31").Value = 3142.8
$ws.Range("L31").Value = 5924
$ws.Range("M31").Value = -2847.8
$ws.Range("N31").Value = -6514
$ws.Range("H34").Value = 4659.8184
$ws.Range("I34").Value = 3142.8
$ws.Range("J34").Value = 5924
$ws.Range("K34").Value = 3142.8
$ws.Range("L34").Value = 5924
$ws.Range("M34").Value = -2940.8
$ws.Range("N34").Value = -6328
$ws.Range("H50").Value = 12749.5
$ws.Range("J50").Value = 12749.5
$ws.Range("L50").Value = 12749.5
$ws.Range("N50").Value = -13999.5
$ws.Range("H51").Value = 9439.700000000001
$ws.Range("I51").Value = 8800
$ws.Range("J51").Value = 9599.625
$ws.Range("K51").Value = 8800
$ws.Range("L51").Value = 9599.625
$ws.Range("M51").Value = -8064
$ws.Range("N51").Value = -11071.625
$ws.Range("H59").Value = 15308.1
$ws.Range("J59").Value = 15231.223
$ws.Range("L59").Value = 15231.223
$ws.Range("N59").Value = -17521.223
$ws.Range("H61").Value = 9439.700000000001
$ws.Range("I61").Value = 8800
$ws.Range("J61").Value = 9599.625
$ws.Range("K61").Value = 8800
$ws.Range("L61").Value = 9599.625
$ws.Range("M61").Value = -8452
$ws.Range("N61").Value = -10295.625
$ws.Range("H74").Value = 14366
$ws.Range("J74").Value = 16406.223
$ws.Range("L74").Value = 16406.223
$ws.Range("N74").Value = -18154.223
$ws.Range("H77").Value = 14366
$ws.Range("J77").Value = 16406.223
$ws.Range("L77").Value = 49218.66900000001
$ws.Range("N77").Value = -57954.66900000001
$ws.Range("H132").Value = 2917.762
$ws.Range("I132").Value = 2601.7222
$ws.Range("J132").Value = 4814
$ws.Range("K132").Value = 7805.1666
$ws.Range("L132").Value = 14442
$ws.Range("M132").Value = -5275.1666
$ws.Range("N132").Value = -19502

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 595.9048
$ws.Range("I113").Value = 619.26086
$ws.Range("J113").Value = 567.6316
$ws.Range("K113").Value = 1857.78258
$ws.Range("L113").Value = 1702.8948
$ws.Range("M113").Value = 312.2174199999999
$ws.Range("N113").Value = -6042.8948
$ws.Range("H122").Value = 7003908.5
$ws.Range("I122").Value = 13889478
$ws.Range("J122").Value = 883402.9399999999
$ws.Range("K122").Value = 125005302
$ws.Range("L122").Value = 7950626.459999999
$ws.Range("M122").Value = -125002852
$ws.Range("N122").Value = -7955526.459999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2323.077
$ws.Range("I102").Value = 1914.2858
$ws.Range("J102").Value = 2800
$ws.Range("K102").Value = 1914.2858
$ws.Range("L102").Value = 2800
$ws.Range("M102").Value = -292.2858000000001
$ws.Range("N102").Value = -6044
$ws.Range("H122").Value = 3623.8333
$ws.Range("I122").Value = 3808.65
$ws.Range("K122").Value = 11425.95
$ws.Range("M122").Value = -8975.950000000001
$ws.Range("H132").Value = 2522.2334
$ws.Range("I132").Value = 2081.2083
$ws.Range("J132").Value = 4286.3335
$ws.Range("K132").Value = 6243.624899999999
$ws.Range("L132").Value = 12859.0005
$ws.Range("M132").Value = -3713.624899999999
$ws.Range("N132").Value = -17919.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8509.700000000001
$ws.Range("I132").Value = 8788.666999999999
$ws.Range("K132").Value = 26366.001
$ws.Range("M132").Value = -23836.001
$ws.Range("H136").Value = 2207.4
$ws.Range("I136").Value = 1175.9166
$ws.Range("J136").Value = 6333.3335
$ws.Range("K136").Value = 3527.7498
$ws.Range("L136").Value = 19000.0005
$ws.Range("M136").Value = -977.7498000000001
$ws.Range("N136").Value = -24100.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3267.838
$ws.Range("I132").Value = 3165.5862
$ws.Range("J132").Value = 3638.5
$ws.Range("K132").Value = 9496.758600000001
$ws.Range("L132").Value = 10915.5
$ws.Range("M132").Value = -6966.758600000001
$ws.Range("N132").Value = -15975.5
